$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing header
# cell onto the new header cells before writing their text.
$headerSrc = $ws.Range("AC1")
$headerSrc.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$lastRow = 51
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 91
    $ws.Cells.Item($row, 31).Value = 71
    $ws.Cells.Item($row, 32).Value = 0
}
